# Applies the "updated alignment & fonts" commit to the poster slide.
#
# Shape index map (1-based, Shapes collection document order) for the
# shapes touched by this edit - resolved by inspecting the underlying
# OOXML of before.pptx:
#   72 -> id=260 "TextBox 259"  ("Dynamic A*" legend label)
#   73 -> id=261 "TextBox 260"  ("Static A*" legend label)
#   89 -> id=9   "TextBox 8"    ("19.1 s" callout, green)
#   90 -> id=166 "TextBox 165"  ("3.4 s" callout, dark blue)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Dynamic A*" legend label: swap font family to Montserrat Semi Bold ---
$dynLabel = $s.Shapes.Item(72)
$dynLabel.TextFrame.TextRange.Font.Name = "Montserrat Semi Bold"

# --- "Static A*" legend label: drop bold + swap font family ---
$staticLabel = $s.Shapes.Item(73)
$staticLabel.TextFrame.TextRange.Font.Bold = $False
$staticLabel.TextFrame.TextRange.Font.Name = "Montserrat Semi Bold"

# --- "19.1 s" callout: re-align/resize box and bump font size to 24pt ---
$greenCallout = $s.Shapes.Item(89)
$greenCallout.Left = 535.553924567874
$greenCallout.Top = 2245.4873228346455
$greenCallout.Width = 95.51668930338583
$greenCallout.Height = 36.351575853149605
$greenCallout.TextFrame.TextRange.Font.Size = 24

# --- "3.4 s" callout: re-align/resize box and bump font size to 24pt ---
$blueCallout = $s.Shapes.Item(90)
$blueCallout.Left = 117.36677170354331
$blueCallout.Top = 2444.2581102362205
$blueCallout.Width = 72.1467704835433
$blueCallout.Height = 36.351575853149605
$blueCallout.TextFrame.TextRange.Font.Size = 24
